$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Insert a new row at position 5 (shifts existing rows 5-13 down to 6-14,
# matching the diff's dimension change from A1:O13 to A1:O14)
$ws.Rows.Item(5).Insert()

# Row 4's D4:K4 were blank placeholder cells; the edit turns them into
# explicit "nan" text (same convention used everywhere else in this sheet)
$ws.Range("D4").Value = "nan"
$ws.Range("E4").Value = "nan"
$ws.Range("F4").Value = "nan"
$ws.Range("G4").Value = "nan"
$ws.Range("H4").Value = "nan"
$ws.Range("I4").Value = "nan"
$ws.Range("J4").Value = "nan"
$ws.Range("K4").Value = "nan"

# Populate the newly inserted row 5 with the new service record (range 151-300).
# A5/B5/C5 hold numeric-looking text ("1"/"151"/"300") that must stay TEXT
# (matching every other row in this column), not be auto-coerced to numbers,
# so they're entered with a leading quote then the formatting is reset to
# plain (no quote-prefix style, no explicit number format) to mirror the
# original file's styling.
$ws.Range("A5").Value = "'1"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "'151"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "'300"
$ws.Range("C5").ClearFormats()

$ws.Range("L5").Value = "29\9\2024"
$ws.Range("M5").Value = "زياره توكيل"
$ws.Range("N5").Value = "اعاده عيار ماكينه"
$ws.Range("O5").Value = "م.صيام"

Write-Output "done"
